# Fruta / hortaliza, semanal
# Insert two new weekly records into the "Macroferia Regional de Talca - Plátano"
# data table, right before the existing row 514, shifting the remaining rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 514-515 (pushes old rows 514:544 down to 516:546).
$ws.Range("A514:A515").EntireRow.Insert()

# New row 514: Plátano "Pintón", week of 2022-02-18 (serial 44610)
$ws.Range("A514").Value = 5
$ws.Range("B514").Value = "Macroferia Regional de Talca"
$ws.Range("C514").Value = "Maule"
$ws.Range("D514").Value = 44610
$ws.Range("E514").Value = 7
$ws.Range("F514").Value = "Fruta"
$ws.Range("G514").Value = 100108
$ws.Range("H514").Value = "Tropicales y subtropicales"
$ws.Range("I514").Value = 100108006
$ws.Range("J514").Value = "Plátano"
$ws.Range("K514").Value = "Sin especificar"
$ws.Range("L514").Value = "Pintón"
$ws.Range("M514").Value = 450
$ws.Range("N514").Value = 16000
$ws.Range("O514").Value = 16000
$ws.Range("P514").Value = 16000
$ws.Range("Q514").Value = "$/caja 20 kilos"
$ws.Range("R514").Value = "Ecuador"
$ws.Range("S514").Value = 800
$ws.Range("T514").Value = 20

# New row 515: Plátano "Primera Pintón", same week (serial 44610)
$ws.Range("A515").Value = 5
$ws.Range("B515").Value = "Macroferia Regional de Talca"
$ws.Range("C515").Value = "Maule"
$ws.Range("D515").Value = 44610
$ws.Range("E515").Value = 7
$ws.Range("F515").Value = "Fruta"
$ws.Range("G515").Value = 100108
$ws.Range("H515").Value = "Tropicales y subtropicales"
$ws.Range("I515").Value = 100108006
$ws.Range("J515").Value = "Plátano"
$ws.Range("K515").Value = "Sin especificar"
$ws.Range("L515").Value = "Primera Pintón"
$ws.Range("M515").Value = 250
$ws.Range("N515").Value = 17000
$ws.Range("O515").Value = 17000
$ws.Range("P515").Value = 17000
$ws.Range("Q515").Value = "$/caja 20 kilos"
$ws.Range("R515").Value = "Ecuador"
$ws.Range("S515").Value = 850
$ws.Range("T515").Value = 20
